# This script applies the PEBCOM workbook update:
#  1. Inserts a new record as row 21 (Caso -204, PARAGUAY /ALT/ 5549 ...),
#     shifting all subsequent rows down by one.
#  2. Appends four new records (Casos 5973, 5989, 5996, 5998) after the
#     former last row, which is now row 54.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell {
    param($Sheet, $Row, $Col, $Text)
    $cell = $Sheet.Cells.Item($Row, $Col)
    $cell.NumberFormat = "@"
    $cell.Value = $Text
}

function Set-RowData {
    param($Sheet, $Row, $Values)
    Set-TextCell $Sheet $Row 1  $Values[0]   # A Caso
    Set-TextCell $Sheet $Row 2  $Values[1]   # B F. De Reclamo
    Set-TextCell $Sheet $Row 3  $Values[2]   # C Direccion
    Set-TextCell $Sheet $Row 4  $Values[3]   # D Comuna
    Set-TextCell $Sheet $Row 5  $Values[4]   # E OT
    Set-TextCell $Sheet $Row 6  $Values[5]   # F Proveedor Asignado
    Set-TextCell $Sheet $Row 7  $Values[6]   # G Estado
    Set-TextCell $Sheet $Row 8  $Values[7]   # H Observaciones
    Set-TextCell $Sheet $Row 9  $Values[8]   # I Attachments
    Set-TextCell $Sheet $Row 10 $Values[9]   # J Tipo de tarea
    Set-TextCell $Sheet $Row 11 $Values[10]  # K Equipo
    Set-TextCell $Sheet $Row 12 $Values[11]  # L Tipo de Elemento
    $Sheet.Cells.Item($Row, 13).Value = $Values[12]  # M Coordenada_X
    $Sheet.Cells.Item($Row, 14).Value = $Values[13]  # N Coordenada_Y
}

# 1) Insert the new row at position 21 (pushes OLLEROS 2952 and everything
#    below it down by one row).
$ws.Rows.Item(21).Insert()

Set-RowData $ws 21 @("-204","12/31/2023","PARAGUAY /ALT/ 5549","106594 - PALERMO","799540519","PEBCOM","Pendiente","Recambio de columna","0","","","",-58.434516,-34.576579)

# 2) Append the four new trailing rows (55-58), right after the old last
#    row (53) which is now row 54.
Set-RowData $ws 55 @("5973","6/4/2025","PALOS 432","4","807168105","PEBCOM","Pendiente","Columna inclinada","1","Cambio","Sin equipos","Pasante",-58.362579,-34.635096)
Set-RowData $ws 56 @("5989","6/4/2025","BONIFACIO, JOSE 2647","7","807168099","PEBCOM","Pendiente","Picada","1","Cambio","Sin equipos","Pasante",-58.464608,-34.633383)
Set-RowData $ws 57 @("5996","6/4/2025","BACACAY 2205","7","807187775","PEBCOM","Pendiente","Picada","1","Cambio","Sin equipos","Pasante",-58.461271,-34.625615)
Set-RowData $ws 58 @("5998","6/4/2025","TRELLES, MANUEL R. 776","7","807187772","PEBCOM","Pendiente","Picada","1","Cambio","Sin equipos","Pasante",-58.459293,-34.617445)
